$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 358.33334
$ws.Range("I2").Value = 287.5
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 287.5
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -174.5
$ws.Range("N2").Value = -726
$ws.Range("H40").Value = 4264.2856
$ws.Range("J40").Value = 4606.25
$ws.Range("L40").Value = 4606.25
$ws.Range("N40").Value = -4956.25
$ws.Range("H107").Value = 784.6667
$ws.Range("I107").Value = 590.05884
$ws.Range("J107").Value = 1611.75
$ws.Range("K107").Value = 590.05884
$ws.Range("L107").Value = 1611.75
$ws.Range("M107").Value = 1329.94116
$ws.Range("N107").Value = -5451.75
$ws.Range("H113").Value = 20837042
$ws.Range("I113").Value = 5004130
$ws.Range("K113").Value = 5004130
$ws.Range("M113").Value = -5000876
$ws.Range("H137").Value = 3588.743
$ws.Range("I137").Value = 756.0769
$ws.Range("J137").Value = 5262.591
$ws.Range("K137").Value = 2268.2307
$ws.Range("L137").Value = 15787.773
$ws.Range("M137").Value = 281.7692999999999
$ws.Range("N137").Value = -20887.773
$ws.Range("H138").Value = 3132.3635
$ws.Range("I138").Value = 2202.5
$ws.Range("J138").Value = 3246.2246
$ws.Range("K138").Value = 6607.5
$ws.Range("L138").Value = 9738.6738
$ws.Range("M138").Value = -1467.5
$ws.Range("N138").Value = -20018.6738

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6107996
$ws.Range("I32").Value = 7252368.5
$ws.Range("K32").Value = 7252368.5
$ws.Range("M32").Value = -7252081.5
$ws.Range("H61").Value = 23868578
$ws.Range("I61").Value = 55560212
$ws.Range("K61").Value = 55560212
$ws.Range("M61").Value = -55560000
$ws.Range("H110").Value = 1722.3846
$ws.Range("I110").Value = 1490.3636
$ws.Range("J110").Value = 2998.5
$ws.Range("K110").Value = 1490.3636
$ws.Range("L110").Value = 2998.5
$ws.Range("M110").Value = 554.6364000000001
$ws.Range("N110").Value = -7088.5
$ws.Range("H132").Value = 5877.3555
$ws.Range("I132").Value = 2773.75
$ws.Range("K132").Value = 8321.25
$ws.Range("M132").Value = -5791.25
$ws.Range("H136").Value = 23868578
$ws.Range("I136").Value = 55560212
$ws.Range("K136").Value = 166680636
$ws.Range("M136").Value = -166678086

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 16332.25
$ws.Range("J99").Value = 4699.6665
$ws.Range("L99").Value = 4699.6665
$ws.Range("N99").Value = -7695.6665
$ws.Range("H105").Value = 1499.5
$ws.Range("I105").Value = 1499.5
$ws.Range("K105").Value = 1499.5
$ws.Range("M105").Value = 247.5
$ws.Range("H134").Value = 29975.195
$ws.Range("I134").Value = 2305.879
$ws.Range("J134").Value = 334337.66
$ws.Range("K134").Value = 6917.637
$ws.Range("L134").Value = 1003012.98
$ws.Range("M134").Value = -4382.637
$ws.Range("N134").Value = -1008082.98

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 6723.4
$ws.Range("J26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("N26").Value = -8574
$ws.Range("H31").Value = 590931.3
$ws.Range("I31").Value = 14716.637
$ws.Range("J31").Value = 809495.5
$ws.Range("K31").Value = 14716.637
$ws.Range("L31").Value = 809495.5
$ws.Range("M31").Value = -14421.637
$ws.Range("N31").Value = -810085.5
$ws.Range("H34").Value = 590931.3
$ws.Range("I34").Value = 14716.637
$ws.Range("J34").Value = 809495.5
$ws.Range("K34").Value = 14716.637
$ws.Range("L34").Value = 809495.5
$ws.Range("M34").Value = -14514.637
$ws.Range("N34").Value = -809899.5
$ws.Range("H58").Value = 888.25
$ws.Range("I58").Value = 599.3333
$ws.Range("J58").Value = 1755
$ws.Range("K58").Value = 599.3333
$ws.Range("L58").Value = 1755
$ws.Range("M58").Value = -396.3333
$ws.Range("N58").Value = -2161
$ws.Range("H86").Value = 3299
$ws.Range("I86").Value = 3156.2856
$ws.Range("K86").Value = 3156.2856
$ws.Range("M86").Value = -2033.2856
$ws.Range("H89").Value = 3299
$ws.Range("I89").Value = 3156.2856
$ws.Range("K89").Value = 15781.428
$ws.Range("M89").Value = -10165.428
$ws.Range("H132").Value = 8791.647000000001
$ws.Range("I132").Value = 4247
$ws.Range("K132").Value = 12741
$ws.Range("M132").Value = -10211
$ws.Range("H136").Value = 888.25
$ws.Range("I136").Value = 599.3333
$ws.Range("J136").Value = 1755
$ws.Range("K136").Value = 1797.9999
$ws.Range("L136").Value = 5265
$ws.Range("M136").Value = 752.0001
$ws.Range("N136").Value = -10365

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 500
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H87").Value = 45007
$ws.Range("I87").Value = 70014
$ws.Range("K87").Value = 210042
$ws.Range("M87").Value = -208794
$ws.Range("H89").Value = 500
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H90").Value = 45007
$ws.Range("I90").Value = 70014
$ws.Range("K90").Value = 630126
$ws.Range("M90").Value = -623886
$ws.Range("H117").Value = 4768.5
$ws.Range("J117").Value = 4768.5
$ws.Range("L117").Value = 14305.5
$ws.Range("N117").Value = -21189.5
$ws.Range("H122").Value = 732.55554
$ws.Range("I122").Value = 574.7778
$ws.Range("J122").Value = 890.3333
$ws.Range("K122").Value = 5173.000199999999
$ws.Range("L122").Value = 8012.9997
$ws.Range("M122").Value = -2723.000199999999
$ws.Range("N122").Value = -12912.9997
$ws.Range("H132").Value = 2261.1052
$ws.Range("I132").Value = 2324.4666
$ws.Range("J132").Value = 2023.5
$ws.Range("K132").Value = 20920.1994
$ws.Range("L132").Value = 18211.5
$ws.Range("M132").Value = -18390.1994
$ws.Range("N132").Value = -23271.5
$ws.Range("H133").Value = 7000
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 276.0625
$ws.Range("H102").Value = 3542.7083
$ws.Range("I102").Value = 3122.2856
$ws.Range("K102").Value = 3122.2856
$ws.Range("M102").Value = -1500.2856
$ws.Range("H126").Value = 5155.4443
$ws.Range("I126").Value = 5466.3335
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 16399.0005
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -13929.0005
$ws.Range("N126").Value = -19940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 204803.8
$ws.Range("I7").Value = 6002
$ws.Range("K7").Value = 6002
$ws.Range("M7").Value = -5890
$ws.Range("H46").Value = 2040.9286
$ws.Range("I46").Value = 2063.6667
$ws.Range("K46").Value = 2063.6667
$ws.Range("M46").Value = -1875.6667
$ws.Range("H55").Value = 90910140
$ws.Range("I55").Value = 125001070
$ws.Range("J55").Value = 1000.3333
$ws.Range("K55").Value = 125001070
$ws.Range("L55").Value = 1000.3333
$ws.Range("M55").Value = -125000897
$ws.Range("N55").Value = -1346.3333
$ws.Range("H63").Value = 112246.664
$ws.Range("J63").Value = 112246.664
$ws.Range("L63").Value = 112246.664
$ws.Range("N63").Value = -113744.664
$ws.Range("H66").Value = 112246.664
$ws.Range("J66").Value = 112246.664
$ws.Range("L66").Value = 336739.992
$ws.Range("N66").Value = -344227.992
$ws.Range("H100").Value = 11463.889
$ws.Range("I100").Value = 16325
$ws.Range("J100").Value = 9033.333000000001
$ws.Range("K100").Value = 16325
$ws.Range("L100").Value = 9033.333000000001
$ws.Range("M100").Value = -15784
$ws.Range("N100").Value = -10115.333
$ws.Range("H122").Value = 6649.25
$ws.Range("I122").Value = 5921.278
$ws.Range("K122").Value = 17763.834
$ws.Range("M122").Value = -15313.834
$ws.Range("H126").Value = 204803.8
$ws.Range("I126").Value = 6002
$ws.Range("K126").Value = 18006
$ws.Range("M126").Value = -15536
$ws.Range("H132").Value = 3972812.8
$ws.Range("I132").Value = 529369.2
$ws.Range("J132").Value = 14303143
$ws.Range("K132").Value = 1588107.6
$ws.Range("L132").Value = 42909429
$ws.Range("M132").Value = -1585577.6
$ws.Range("N132").Value = -42914489

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 614.775
$ws.Range("I113").Value = 539.38464
$ws.Range("J113").Value = 754.7857
$ws.Range("K113").Value = 1618.15392
$ws.Range("L113").Value = 2264.3571
$ws.Range("M113").Value = 551.84608
$ws.Range("N113").Value = -6604.3571
$ws.Range("H122").Value = 5141.85
$ws.Range("I122").Value = 4048.9412
$ws.Range("K122").Value = 12146.8236
$ws.Range("M122").Value = -9696.8236
$ws.Range("H126").Value = 4941.923
$ws.Range("I126").Value = 3080.524
$ws.Range("J126").Value = 12759.8
$ws.Range("K126").Value = 9241.572
$ws.Range("L126").Value = 38279.39999999999
$ws.Range("M126").Value = -6771.572
$ws.Range("N126").Value = -43219.39999999999
$ws.Range("H132").Value = 296679.47
$ws.Range("I132").Value = 2059.9333
$ws.Range("J132").Value = 2506326
$ws.Range("K132").Value = 6179.7999
$ws.Range("L132").Value = 7518978
$ws.Range("M132").Value = -3649.7999
$ws.Range("N132").Value = -7524038
